$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("overview")

# Update the "Structural balance demo" link in E5 to point to the fa2025 materials
$ws.Range("E5").Value = "[Structural balance demo](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2025&branch=main&urlpath=lab%2Ftree%2Fdemog180-fa2025%2Flecture%2F20250908_structural_balance%2Fstructural_balance_in_the_small_slashdot_network.ipynb)"

# Move "Mini Project 02: Complete network data" from G7 down to G8
$ws.Range("G8").Value = "Mini Project 02: Complete network data"
$ws.Range("G7").Value = $null
